$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'62.951.89"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -5.91%  '
$ws.Range('D3').Value = "'3.304.11"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -6.30%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = "'546.74"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.58%  '
$ws.Range('D6').Value = "'169.25"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -9.86%  '
$ws.Range('D7').Value = "'0.604"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.80%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = "'3.296.15"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -6.41%  '
$ws.Range('D10').Value = "'0.608"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.05%  '
$ws.Range('E11').Value = '  -6.15%  '
$ws.Range('D12').Value = "'53.61"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.19%  '
$ws.Range('E13').Value = '  -4.94%  '
$ws.Range('D14').Value = "'8.80"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -5.84%  '
$ws.Range('D15').Value = "'3.839.48"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -6.19%  '
$ws.Range('D17').Value = "'3.308.42"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -6.40%  '
$ws.Range('D18').Value = "'17.45"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -6.34%  '
$ws.Range('D19').Value = "'62.936.89"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -6.01%  '
$ws.Range('E20').Value = '  -5.41%  '
$ws.Range('D21').Value = "'0.959"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.92%  '
$ws.Range('D22').Value = "'398.72"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.30%  '
$ws.Range('D23').Value = "'3.98"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.55%  '
$ws.Range('D24').Value = "'4.21"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.02%  '
$ws.Range('D25').Value = "'13.02"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +6.09%  '
$ws.Range('D26').Value = "'81.47"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -5.19%  '
$ws.Range('D27').Value = "'10.61"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.86%  '
$ws.Range('E28').Value = '  -7.22%  '
$ws.Range('D29').Value = "'8.53"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -6.18%  '
$ws.Range('D30').Value = "'28.76"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.32%  '
$ws.Range('D31').Value = "'6.49"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.25%  '
$ws.Range('D32').Value = "'576.41"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -9.47%  '
$ws.Range('D33').Value = "'11.15"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.91%  '
$ws.Range('E34').Value = '  -6.42%  '
$ws.Range('D35').Value = "'57.47"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.68%  '
$ws.Range('E36').Value = '  -0.53%  '
$ws.Range('D37').Value = "'1.00"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.17%  '
$ws.Range('D38').Value = "'35.36"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -7.64%  '
$ws.Range('D39').Value = "'3.43"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.59%  '
$ws.Range('B40').Value = 'TheGraph'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D40').Value = "'0.362"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -6.39%  '
$ws.Range('B41').Value = 'PEPE'
$ws.Range('C41').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D41').Value = "'0.0₃0723"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -12.45%  '
$ws.Range('D42').Value = "'3.090.44"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.12%  '
$ws.Range('D43').Value = "'0.999"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.17%  '
$ws.Range('D44').Value = "'2.76"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.32%  '
$ws.Range('D45').Value = "'2.44"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -7.27%  '
$ws.Range('D46').Value = "'3.16"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -6.30%  '
$ws.Range('E47').Value = '  -5.27%  '
$ws.Range('D48').Value = "'2.57"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -6.41%  '
$ws.Range('E49').Value = '  -4.89%  '
$ws.Range('D50').Value = "'130.92"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -6.02%  '
$ws.Range('D51').Value = "'7.93"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -6.93%  '
